$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) First three rows: "100" -> "0M", "0" -> "0M", "146" -> "0M"
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# 2) Insert 10 new rows after row 3 (before the row that was originally row 4, text "0")
$newValues = @("102", "0.00002", "0.00005", "0.00004", "0.00001", "0.00003", "0.00004", "0.00004", "0.00389", "100.0")
$insertPos = 4
foreach ($val in $newValues) {
    $anchorRow = $t.Rows.Item($insertPos)
    $newRow = $t.Rows.Add($anchorRow)
    $newRow.Cells.Item(1).Range.Text = $val
    $insertPos = $insertPos + 1
}

# 3) Near the end of the table, the row that had "100" followed by many
#    tab-separated values collapses down to just "100".
$n = $t.Rows.Count
$row1 = $t.Rows.Item($n - 2)
$row1.Cells.Item(1).Range.Text = "100"

# 4) The following row (originally "2" followed by many tab-separated
#    values) collapses down to just "0".
$row2 = $t.Rows.Item($n - 1)
$row2.Cells.Item(1).Range.Text = "0"

# 5) The last (previously empty) row gets the text "146".
$row3 = $t.Rows.Item($n)
$row3.Cells.Item(1).Range.Text = "146"
